$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one
# day (46074 -> 46075) for every data row (rows 2 through 489).
$ws.Range("C2:C489").Value = 46075
